try {
    $wb = $excel.ActiveWorkbook

    $wsVertical = $wb.Worksheets.Item("ValidExampleVerticalData")
    $wsSetting  = $wb.Worksheets.Item("ValidSetting")

    # ---------------------------------------------------------------
    # ValidSetting ("#active" column gets added, like ValidExampleVerticalData already has)
    # ---------------------------------------------------------------

    # Insert a new header/metadata row above the existing data (old row 2 -> row 3, etc.)
    $wsSetting.Rows.Item(2).Insert()

    # Insert a new blank column D (old D -> E, numeric settings column)
    $wsSetting.Columns.Item(4).Insert()

    # The freshly inserted column inherited formatting from its left neighbour; reset the
    # data-row cells in the new column back to an unformatted/default style (matches the
    # target which has no explicit style override there) by pasting formats from a cell
    # that has never been touched.
    $blankFormat = $wsSetting.Cells.Item(500, 500)
    $blankFormat.Copy()
    $wsSetting.Range($wsSetting.Cells.Item(3, 4), $wsSetting.Cells.Item(23, 4)).PasteSpecial(-4122)

    # Fill in the new "#active" row (row 2): id column blank, then #active/bool/skip/true
    $wsSetting.Cells.Item(2, 2).Value = "#active"
    $wsSetting.Cells.Item(2, 3).Value = "bool"
    $wsSetting.Cells.Item(2, 4).Value = "skip"

    # "true" would auto-coerce to a boolean if typed directly, so compute it as a text
    # formula first and then collapse it down to a plain cached value.
    $wsSetting.Cells.Item(2, 5).Formula = '="true"'
    $wsSetting.Cells.Item(2, 5).Copy()
    $wsSetting.Cells.Item(2, 5).PasteSpecial(-4163)

    # Fill the new column D for every data row with "both"
    for ($r = 3; $r -le 23; $r++) {
        $wsSetting.Cells.Item($r, 4).Value = "both"
    }

    # ---------------------------------------------------------------
    # Selection / active-sheet bookkeeping
    # ---------------------------------------------------------------

    $wsVertical.Activate() | Out-Null
    $wsVertical.Range("B2").Select() | Out-Null

    $wsSetting.Activate() | Out-Null
    $wsSetting.Range("G19").Select() | Out-Null
} catch {
    Write-Host "ERROR: $_"
    throw
}
